$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 9316.333000000001
$ws.Range("I31").Value = 1398
$ws.Range("J31").Value = 10306.125
$ws.Range("K31").Value = 4194
$ws.Range("L31").Value = 30918.375
$ws.Range("M31").Value = -3964
$ws.Range("N31").Value = -31378.375
# Row 74
$ws.Range("H74").Value = 21723
$ws.Range("I74").Value = 9321
$ws.Range("J74").Value = 34125
$ws.Range("K74").Value = 9321
$ws.Range("L74").Value = 34125
$ws.Range("M74").Value = -8385
$ws.Range("N74").Value = -35997
# Row 77
$ws.Range("H77").Value = 21723
$ws.Range("I77").Value = 9321
$ws.Range("J77").Value = 34125
$ws.Range("K77").Value = 46605
$ws.Range("L77").Value = 170625
$ws.Range("M77").Value = -41925
$ws.Range("N77").Value = -179985
# Row 116
$ws.Range("H116").Value = 10734.952
$ws.Range("I116").Value = 3639.2
$ws.Range("J116").Value = 12952.375
$ws.Range("K116").Value = 3639.2
$ws.Range("L116").Value = 12952.375
$ws.Range("M116").Value = -197.1999999999998
$ws.Range("N116").Value = -19836.375
# Row 137
$ws.Range("H137").Value = 1646.2941
$ws.Range("I137").Value = 1320.5
$ws.Range("K137").Value = 3961.5
$ws.Range("M137").Value = -1411.5
# Row 138
$ws.Range("H138").Value = 2419.4211
$ws.Range("I138").Value = 2130.3333
$ws.Range("J138").Value = 3228.8667
$ws.Range("K138").Value = 6390.999899999999
$ws.Range("L138").Value = 9686.6001
$ws.Range("M138").Value = -1250.999899999999
$ws.Range("N138").Value = -19966.6001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 490
$ws.Range("I26").Value = 490
$ws.Range("K26").Value = 490
$ws.Range("M26").Value = -160
# Row 45
$ws.Range("H45").Value = 4197.4
$ws.Range("I45").Value = 3929.6667
$ws.Range("K45").Value = 3929.6667
$ws.Range("M45").Value = -3552.6667
# Row 61
$ws.Range("H61").Value = 5602.4
$ws.Range("J61").Value = 6670.6665
$ws.Range("L61").Value = 6670.6665
$ws.Range("N61").Value = -7094.6665
# Row 63
$ws.Range("H63").Value = 6220.4546
$ws.Range("I63").Value = 5936.3335
$ws.Range("K63").Value = 5936.3335
$ws.Range("M63").Value = -5250.3335
# Row 66
$ws.Range("H66").Value = 6220.4546
$ws.Range("I66").Value = 5936.3335
$ws.Range("K66").Value = 29681.6675
$ws.Range("M66").Value = -26249.6675
# Row 125
$ws.Range("H125").Value = 79995
$ws.Range("J125").Value = 79995
$ws.Range("L125").Value = 79995
$ws.Range("N125").Value = -89835
# Row 132
$ws.Range("H132").Value = 17235.926
$ws.Range("I132").Value = 21277.076
$ws.Range("J132").Value = 1937.2858
$ws.Range("K132").Value = 63831.228
$ws.Range("L132").Value = 5811.857400000001
$ws.Range("M132").Value = -61301.228
$ws.Range("N132").Value = -10871.8574
# Row 136
$ws.Range("H136").Value = 5602.4
$ws.Range("J136").Value = 6670.6665
$ws.Range("L136").Value = 20011.9995
$ws.Range("N136").Value = -25111.9995

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 22562.25
$ws.Range("I7").Value = 151
$ws.Range("J7").Value = 30032.666
$ws.Range("K7").Value = 151
$ws.Range("L7").Value = 30032.666
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -30258.666
# Row 22
$ws.Range("H22").Value = 4816.0557
$ws.Range("J22").Value = 6697
$ws.Range("L22").Value = 6697
$ws.Range("N22").Value = -7043

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 4071.3845
$ws.Range("I12").Value = 3259.4285
$ws.Range("J12").Value = 5018.6665
$ws.Range("K12").Value = 3259.4285
$ws.Range("L12").Value = 5018.6665
$ws.Range("M12").Value = -3089.4285
$ws.Range("N12").Value = -5358.6665
# Row 31
$ws.Range("H31").Value = 3578.8696
$ws.Range("I31").Value = 1681.909
$ws.Range("J31").Value = 5317.75
$ws.Range("K31").Value = 1681.909
$ws.Range("L31").Value = 5317.75
$ws.Range("M31").Value = -1386.909
$ws.Range("N31").Value = -5907.75
# Row 34
$ws.Range("H34").Value = 3578.8696
$ws.Range("I34").Value = 1681.909
$ws.Range("J34").Value = 5317.75
$ws.Range("K34").Value = 1681.909
$ws.Range("L34").Value = 5317.75
$ws.Range("M34").Value = -1479.909
$ws.Range("N34").Value = -5721.75
# Row 103
$ws.Range("H103").Value = 54208.2
$ws.Range("I103").Value = 50260.5
$ws.Range("J103").Value = 69999
$ws.Range("K103").Value = 50260.5
$ws.Range("L103").Value = 69999
$ws.Range("M103").Value = -49088.5
$ws.Range("N103").Value = -72343

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 1500
$ws.Range("J33").Value = 2500
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15566
# Row 52
$ws.Range("H52").Value = 900
$ws.Range("J52").Value = 900
$ws.Range("L52").Value = 2700
$ws.Range("N52").Value = -3232
# Row 80
$ws.Range("H80").Value = 4999.8
$ws.Range("J80").Value = 5250
$ws.Range("L80").Value = 15750
$ws.Range("N80").Value = -17622
# Row 83
$ws.Range("H83").Value = 4999.8
$ws.Range("J83").Value = 5250
$ws.Range("L83").Value = 47250
$ws.Range("N83").Value = -56610
# Row 118
$ws.Range("H118").Value = 2799.5
$ws.Range("I118").Value = 2799.5
$ws.Range("K118").Value = 8398.5
$ws.Range("M118").Value = -7155.5
# Row 119
$ws.Range("H119").Value = 4286
$ws.Range("I119").Value = 4286
$ws.Range("K119").Value = 12858
$ws.Range("M119").Value = -8020
# Row 125
$ws.Range("H125").Value = 28333.334
$ws.Range("I125").Value = 30000
$ws.Range("J125").Value = 27500
$ws.Range("K125").Value = 90000
$ws.Range("L125").Value = 82500
$ws.Range("M125").Value = -85080
$ws.Range("N125").Value = -92340
# Row 130
$ws.Range("H130").Value = 6173.25
$ws.Range("I130").Value = 3757.6
$ws.Range("J130").Value = 10199.333
$ws.Range("K130").Value = 11272.8
$ws.Range("L130").Value = 30597.999
$ws.Range("M130").Value = -6252.799999999999
$ws.Range("N130").Value = -40637.999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 125507.625
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("N107").Value = -5840
# Row 122
$ws.Range("H122").Value = 2662.85
$ws.Range("I122").Value = 2280.1765
$ws.Range("J122").Value = 4831.3335
$ws.Range("K122").Value = 6840.529500000001
$ws.Range("L122").Value = 14494.0005
$ws.Range("M122").Value = -4390.529500000001
$ws.Range("N122").Value = -19394.0005
# Row 132
$ws.Range("H132").Value = 49872.668
$ws.Range("I132").Value = 57629.445
$ws.Range("J132").Value = 3332
$ws.Range("K132").Value = 172888.335
$ws.Range("L132").Value = 9996
$ws.Range("M132").Value = -170358.335
$ws.Range("N132").Value = -15056

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 2916.5
$ws.Range("I113").Value = 1473.6666
$ws.Range("J113").Value = 5080.75
$ws.Range("K113").Value = 4420.9998
$ws.Range("L113").Value = 15242.25
$ws.Range("M113").Value = -2250.9998
$ws.Range("N113").Value = -19582.25
# Row 124
$ws.Range("H124").Value = 70162.78
$ws.Range("J124").Value = 70162.78
$ws.Range("L124").Value = 70162.78
$ws.Range("N124").Value = -79982.78
# Row 126
$ws.Range("H126").Value = 52384.7
$ws.Range("I126").Value = 64999.062
$ws.Range("J126").Value = 1927.25
$ws.Range("K126").Value = 194997.186
$ws.Range("L126").Value = 5781.75
$ws.Range("M126").Value = -192527.186
$ws.Range("N126").Value = -10721.75
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
